$wb = $excel.ActiveWorkbook

# ---- Sheet "data" (sheet1): add new date column AE (col 31) "25. 5. 2021" ----
$ws1 = $wb.Worksheets.Item("data")
$ws1.Cells.Item(1, 30).Copy()
$ws1.Cells.Item(1, 31).PasteSpecial(-4122)
$ws1.Cells.Item(1, 31).Value = "25. 5. 2021"

$ae1 = @{
  2 = 0.52
  3 = 0.27
  4 = 0.21
  5 = 0.69
  6 = 0.2
  7 = 0.11
  8 = 0.56
  9 = 0.28
  10 = 0.16
  11 = 0.38
  12 = 0.31
  13 = 0.31
  14 = 0.38
  15 = 0.25
  16 = 0.37
  17 = 0.47
  18 = 0.33
  19 = 0.2
  20 = 0.6
  21 = 0.26
  22 = 0.14
  23 = 0.7
  24 = 0.18
  25 = 0.12
  26 = 0.52
  27 = 0.28
  28 = 0.2
  29 = 0.36
  30 = 0.33
  31 = 0.31
  32 = 0.56
  33 = 0.26
  34 = 0.18
  35 = 0.46
  36 = 0.29
  37 = 0.25
  38 = 0.54
  39 = 0.26
  40 = 0.2
  41 = 0.55
  42 = 0.25
  43 = 0.2
  44 = 0.49
  45 = 0.29
  46 = 0.22
  47 = 0.54
  48 = 0.25
  49 = 0.21
  50 = 0.51
  51 = 0.27
  52 = 0.22
  53 = 0.49
  54 = 0.31
  55 = 0.2
  56 = 0.46
  57 = 0.34
  58 = 0.2
}
foreach ($r in $ae1.Keys) {
  $ws1.Cells.Item([int]$r, 31).Value = $ae1[$r]
}

# ---- Sheet "data": corrected AD (col 30) values ----
$adCorr1 = @{
  14 = 0.32
  16 = 0.36
  53 = 0.43
  55 = 0.23
  57 = 0.32
  58 = 0.24
}
foreach ($r in $adCorr1.Keys) {
  $ws1.Cells.Item([int]$r, 30).Value = $adCorr1[$r]
}

# ---- Sheet "data": update title cell A59 ----
$ws1.Cells.Item(59, 1).Value = "Život během pandemie, Počet protektivních aktivit, % respondentů celkově a ve skupinách, aktualizace 1. 6. 2021"

# ---- Sheet "pocetR" (sheet2): add new date column AD (col 30) "25. 5. 2021" ----
$ws2 = $wb.Worksheets.Item("pocetR")
$ws2.Cells.Item(1, 29).Copy()
$ws2.Cells.Item(1, 30).PasteSpecial(-4122)
$ws2.Cells.Item(1, 30).Value = "25. 5. 2021"

$ad2 = @{
  2 = 1975
  3 = 471
  4 = 729
  5 = 775
  6 = 373
  7 = 643
  8 = 719
  9 = 512
  10 = 844
  11 = 619
  12 = 443
  13 = 528
  14 = 1004
  15 = 956
  16 = 1019
  17 = 1029
  18 = 456
  19 = 230
  20 = 260
}
foreach ($r in $ad2.Keys) {
  $ws2.Cells.Item([int]$r, 30).Value = $ad2[$r]
}

# ---- Sheet "pocetR": empty string cell AD21 (new, matches existing row pattern) ----
$ws2.Cells.Item(21, 30).Value = ""

# ---- Sheet "pocetR": corrected AB/AC values ----
$corr2 = @(
  @{Row=2; Col="AC"; Val=2029}
  @{Row=3; Col="AC"; Val=480}
  @{Row=4; Col="AC"; Val=749}
  @{Row=5; Col="AC"; Val=800}
  @{Row=6; Col="AC"; Val=451}
  @{Row=7; Col="AC"; Val=689}
  @{Row=8; Col="AB"; Val=518}
  @{Row=8; Col="AC"; Val=678}
  @{Row=11; Col="AC"; Val=737}
  @{Row=12; Col="AC"; Val=493}
  @{Row=13; Col="AC"; Val=591}
  @{Row=15; Col="AC"; Val=980}
  @{Row=16; Col="AC"; Val=1049}
  @{Row=17; Col="AC"; Val=1058}
  @{Row=18; Col="AC"; Val=466}
  @{Row=20; Col="AC"; Val=267}
)
foreach ($item in $corr2) {
  $colNum = if ($item.Col -eq "AB") { 28 } else { 29 }
  $ws2.Cells.Item($item.Row, $colNum).Value = $item.Val
}

# ---- Sheet "pocetR": update title cell A21 ----
$ws2.Cells.Item(21, 1).Value = "Život během pandemie, Počet protektivních aktivit, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 6. 2021"